# Refresh the crypto price/volume table (Sheet1) with the latest scraped
# values. Price-looking numeric strings are entered with a leading
# apostrophe so Excel keeps them as text (preserving formats such as
# "60.380.37" or a trailing zero like "0.0000140") instead of silently
# reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.380.37"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "2.595.08"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'564.96"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'142.52"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "2.617.64"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "'6.56"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "'0.369"
$ws.Range("E12").Value = "  +7.16%  "
$ws.Range("D13").Value = "'0.152"
$ws.Range("E13").Value = "  -4.65%  "
$ws.Range("D14").Value = "3.055.37"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "60.421.06"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "'23.29"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "2.606.50"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  +8.17%  "
$ws.Range("D20").Value = "'4.64"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'346.42"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "'6.96"
$ws.Range("E22").Value = "  +9.63%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +16.35%  "
$ws.Range("D25").Value = "'63.32"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'7.65"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "0.0₃0784"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "'1.79"
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'6.31"
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("D33").Value = "'161.02"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "'19.42"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").Value = "'4.23"
$ws.Range("E35").Value = "  +5.22%  "
$ws.Range("D36").Value = "'0.958"
$ws.Range("E36").Value = "  +9.27%  "
$ws.Range("E37").Value = "  +4.72%  "
$ws.Range("E38").Value = "  +8.20%  "
$ws.Range("D39").Value = "'37.71"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "'0.859"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.81"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").Value = "'294.09"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'136.92"
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "'0.0977"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'0.0546"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "'19.45"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "'0.0239"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'19.66"
$ws.Range("E51").Value = "  +6.57%  "
